$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2 through 89 from 45208 to 45212
for ($r = 2; $r -le 89; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}

# Update hyperlink formulas in row 2 to point to the renamed files
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/artfynd/A 33491-2023 artfynd.xlsx", "A 33491-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/kartor/A 33491-2023 karta.png", "A 33491-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomål/A 33491-2023 fsc-klagomål.docx", "A 33491-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomålsmail/A 33491-2023 fsc-klagomål mail.docx", "A 33491-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/tillsyn/A 33491-2023 tillsynsbegäran.docx", "A 33491-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/ti,llsynsmail/A 33491-2023 tillsynsbegäran mail.docx", "A 33491-2023")'
